$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string (A1)
$ws.Range('A1').Value = 'Datos actualizados a 19 de Junio de 2020 a las 22:09'

# Update country data rows whose rank / values changed due to new case counts
# Row 4: 'Estados Unidos' (values updated)
$ws.Range('B4').Value = 2286093
$ws.Range('C4').Value = 22442
$ws.Range('D4').Value = 932633
$ws.Range('E4').Value = 1232239
$ws.Range('F4').Value = 0
$ws.Range('G4').Value = 533
$ws.Range('H4').Value = 121221

# Row 16: 'Pakistan' -> 'Mexico'
$ws.Range('A16').Value = 'Mexico'
$ws.Range('B16').Value = 165455
$ws.Range('C16').Value = 5662
$ws.Range('D16').Value = 123095
$ws.Range('E16').Value = 22613
$ws.Range('F16').Value = 0
$ws.Range('G16').Value = 667
$ws.Range('H16').Value = 19747

# Row 17: 'Mexico' -> 'Pakistan'
$ws.Range('A17').Value = 'Pakistan'
$ws.Range('B17').Value = 165062
$ws.Range('C17').Value = 4944
$ws.Range('D17').Value = 61383
$ws.Range('E17').Value = 100450
$ws.Range('F17').Value = 0
$ws.Range('G17').Value = 136
$ws.Range('H17').Value = 3229

# Row 18: 'Francia' (values updated)
$ws.Range('B18').Value = 159452
$ws.Range('C18').Value = 811
$ws.Range('D18').Value = 74117
$ws.Range('E18').Value = 55718
$ws.Range('F18').Value = 0
$ws.Range('G18').Value = 14
$ws.Range('H18').Value = 29617

# Row 22: 'Catar' -> 'Sudafrica'
$ws.Range('A22').Value = 'Sudafrica'
$ws.Range('B22').Value = 87715
$ws.Range('C22').Value = 3825
$ws.Range('D22').Value = 47825
$ws.Range('E22').Value = 38059
$ws.Range('F22').Value = 0
$ws.Range('G22').Value = 94
$ws.Range('H22').Value = 1831

# Row 23: 'Sudafrica' -> 'Catar'
$ws.Range('A23').Value = 'Catar'
$ws.Range('B23').Value = 85462
$ws.Range('C23').Value = 1021
$ws.Range('D23').Value = 65409
$ws.Range('E23').Value = 19960
$ws.Range('F23').Value = 0
$ws.Range('G23').Value = 7
$ws.Range('H23').Value = 93

# Row 75: 'Costa de Marfil' (values updated)
$ws.Range('B75').Value = 6874
$ws.Range('C75').Value = 430
$ws.Range('D75').Value = 2942
$ws.Range('E75').Value = 3883
$ws.Range('F75').Value = 0
$ws.Range('G75').Value = 0
$ws.Range('H75').Value = 49

# Row 103: 'Estonia' -> 'Costa Rica'
$ws.Range('A103').Value = 'Costa Rica'
$ws.Range('B103').Value = 2058
$ws.Range('C103').Value = 119
$ws.Range('D103').Value = 982
$ws.Range('E103').Value = 1064
$ws.Range('F103').Value = 0
$ws.Range('G103').Value = 0
$ws.Range('H103').Value = 12

# Row 104: 'Guayana Francesa' -> 'Estonia'
$ws.Range('A104').Value = 'Estonia'
$ws.Range('B104').Value = 1979
$ws.Range('C104').Value = 2
$ws.Range('D104').Value = 1755
$ws.Range('E104').Value = 155
$ws.Range('F104').Value = 0
$ws.Range('G104').Value = 0
$ws.Range('H104').Value = 69

# Row 105: 'Sri Lanka' -> 'Guayana Francesa'
$ws.Range('A105').Value = 'Guayana Francesa'
$ws.Range('B105').Value = 1969
$ws.Range('C105').Value = 211
$ws.Range('D105').Value = 840
$ws.Range('E105').Value = 1124
$ws.Range('F105').Value = 0
$ws.Range('G105').Value = 0
$ws.Range('H105').Value = 5

# Row 106: 'Costa Rica' -> 'Sri Lanka'
$ws.Range('A106').Value = 'Sri Lanka'
$ws.Range('B106').Value = 1950
$ws.Range('C106').Value = 4
$ws.Range('D106').Value = 1446
$ws.Range('E106').Value = 493
$ws.Range('F106').Value = 0
$ws.Range('G106').Value = 0
$ws.Range('H106').Value = 11

# Row 120: 'Zambia' (values updated)
$ws.Range('B120').Value = 1430
$ws.Range('C120').Value = 14
$ws.Range('D120').Value = 1194
$ws.Range('E120').Value = 225
$ws.Range('F120').Value = 0
$ws.Range('G120').Value = 0
$ws.Range('H120').Value = 11

# Row 121: 'Paraguay' (values updated)
$ws.Range('B121').Value = 1336
$ws.Range('C121').Value = 6
$ws.Range('D121').Value = 741
$ws.Range('E121').Value = 582
$ws.Range('F121').Value = 0
$ws.Range('G121').Value = 0
$ws.Range('H121').Value = 13

# Row 129: 'Yemen' (values updated)
$ws.Range('B129').Value = 919
$ws.Range('C129').Value = 10
$ws.Range('D129').Value = 288
$ws.Range('E129').Value = 380
$ws.Range('F129').Value = 0
$ws.Range('G129').Value = 3
$ws.Range('H129').Value = 251

# Row 141: 'Estado de Palestina' (values updated)
$ws.Range('B141').Value = 675
$ws.Range('C141').Value = 75
$ws.Range('D141').Value = 437
$ws.Range('E141').Value = 235
$ws.Range('F141').Value = 0
$ws.Range('G141').Value = 0
$ws.Range('H141').Value = 3

# Row 144: 'Benin' -> 'Ruanda'
$ws.Range('A144').Value = 'Ruanda'
$ws.Range('B144').Value = 661
$ws.Range('C144').Value = 15
$ws.Range('D144').Value = 351
$ws.Range('E144').Value = 308
$ws.Range('F144').Value = 0
$ws.Range('G144').Value = 0
$ws.Range('H144').Value = 2

# Row 145: 'Ruanda' -> 'Benin'
$ws.Range('A145').Value = 'Benin'
$ws.Range('B145').Value = 650
$ws.Range('C145').Value = 53
$ws.Range('D145').Value = 247
$ws.Range('E145').Value = 392
$ws.Range('F145').Value = 0
$ws.Range('G145').Value = 0
$ws.Range('H145').Value = 11

# Row 147: 'Malaui' -> 'Suazilandia'
$ws.Range('A147').Value = 'Suazilandia'
$ws.Range('B147').Value = 623
$ws.Range('C147').Value = 37
$ws.Range('D147').Value = 276
$ws.Range('E147').Value = 343
$ws.Range('F147').Value = 0
$ws.Range('G147').Value = 0
$ws.Range('H147').Value = 4

# Row 148: 'Suazilandia' -> 'Malaui'
$ws.Range('A148').Value = 'Malaui'
$ws.Range('B148').Value = 592
$ws.Range('C148').Value = 0
$ws.Range('D148').Value = 74
$ws.Range('E148').Value = 510
$ws.Range('F148').Value = 0
$ws.Range('G148').Value = 0
$ws.Range('H148').Value = 8

# Row 150: 'Togo' (values updated)
$ws.Range('B150').Value = 555
$ws.Range('C150').Value = 8
$ws.Range('D150').Value = 361
$ws.Range('E150').Value = 181
$ws.Range('F150').Value = 0
$ws.Range('G150').Value = 0
$ws.Range('H150').Value = 13

# Row 153: 'Reunion' (values updated)
$ws.Range('B153').Value = 504
$ws.Range('C153').Value = 2
$ws.Range('D153').Value = 460
$ws.Range('E153').Value = 43
$ws.Range('F153').Value = 0
$ws.Range('G153').Value = 0
$ws.Range('H153').Value = 1

# Row 160: 'Birmania' (values updated)
$ws.Range('B160').Value = 286
$ws.Range('C160').Value = 23
$ws.Range('D160').Value = 192
$ws.Range('E160').Value = 88
$ws.Range('F160').Value = 0
$ws.Range('G160').Value = 0
$ws.Range('H160').Value = 6

# Row 170: 'Guadalupe' -> 'Angola'
$ws.Range('A170').Value = 'Angola'
$ws.Range('B170').Value = 172
$ws.Range('C170').Value = 6
$ws.Range('D170').Value = 66
$ws.Range('E170').Value = 98
$ws.Range('F170').Value = 0
$ws.Range('G170').Value = 0
$ws.Range('H170').Value = 8

# Row 171: 'Angola' -> 'Guadalupe'
$ws.Range('A171').Value = 'Guadalupe'
$ws.Range('B171').Value = 171
$ws.Range('C171').Value = 0
$ws.Range('D171').Value = 157
$ws.Range('E171').Value = 0
$ws.Range('F171').Value = 0
$ws.Range('G171').Value = 0
$ws.Range('H171').Value = 14

# Row 202: 'Dominica' -> 'Fiyi'
$ws.Range('A202').Value = 'Fiyi'
$ws.Range('B202').Value = 18
$ws.Range('C202').Value = 0
$ws.Range('D202').Value = 18
$ws.Range('E202').Value = 0
$ws.Range('F202').Value = 0
$ws.Range('G202').Value = 0
$ws.Range('H202').Value = 0

# Row 203: 'Fiyi' -> 'Dominica'
$ws.Range('A203').Value = 'Dominica'
$ws.Range('B203').Value = 18
$ws.Range('C203').Value = 0
$ws.Range('D203').Value = 18
$ws.Range('E203').Value = 0
$ws.Range('F203').Value = 0
$ws.Range('G203').Value = 0
$ws.Range('H203').Value = 0

# Row 208: 'Islas Turcas y Caicos' -> 'Santa Sede'
$ws.Range('A208').Value = 'Santa Sede'
$ws.Range('B208').Value = 12
$ws.Range('C208').Value = 0
$ws.Range('D208').Value = 12
$ws.Range('E208').Value = 0
$ws.Range('F208').Value = 0
$ws.Range('G208').Value = 0
$ws.Range('H208').Value = 0

# Row 209: 'Santa Sede' -> 'Islas Turcas y Caicos'
$ws.Range('A209').Value = 'Islas Turcas y Caicos'
$ws.Range('B209').Value = 12
$ws.Range('C209').Value = 0
$ws.Range('D209').Value = 11
$ws.Range('E209').Value = 0
$ws.Range('F209').Value = 0
$ws.Range('G209').Value = 0
$ws.Range('H209').Value = 1
